$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I6").Value = "sd"
$ws.Range("J6").Value = "Statement-non-opinion"
$ws.Range("I16").Value = "sd"
$ws.Range("J16").Value = "Statement-non-opinion"
$ws.Range("I30").Value = "b"
$ws.Range("J30").Value = "Acknowledge (Backchannel)"
$ws.Range("I38").Value = "sv"
$ws.Range("J38").Value = "Statement-opinion"
$ws.Range("I41").Value = "aa"
$ws.Range("J41").Value = "Agree/Accept"
$ws.Range("I50").Value = "%"
$ws.Range("J50").Value = "Uninterpretable"
$ws.Range("I51").Value = "aa"
$ws.Range("J51").Value = "Agree/Accept"
$ws.Range("I55").Value = "aa"
$ws.Range("J55").Value = "Agree/Accept"
$ws.Range("I65").Value = "aa"
$ws.Range("J65").Value = "Agree/Accept"
$ws.Range("I68").Value = "sd"
$ws.Range("J68").Value = "Statement-non-opinion"
$ws.Range("I82").Value = "sd"
$ws.Range("J82").Value = "Statement-non-opinion"
$ws.Range("I85").Value = "qy"
$ws.Range("J85").Value = "Yes-No-Question"
$ws.Range("I91").Value = "sd"
$ws.Range("J91").Value = "Statement-non-opinion"
$ws.Range("I93").Value = "aa"
$ws.Range("J93").Value = "Agree/Accept"
$ws.Range("I107").Value = "sd"
$ws.Range("J107").Value = "Statement-non-opinion"
$ws.Range("I114").Value = "sd"
$ws.Range("J114").Value = "Statement-non-opinion"
$ws.Range("I162").Value = "aa"
$ws.Range("J162").Value = "Agree/Accept"
$ws.Range("I184").Value = "sv"
$ws.Range("J184").Value = "Statement-opinion"
$ws.Range("I198").Value = "ba"
$ws.Range("J198").Value = "Appreciation"
$ws.Range("I224").Value = "aa"
$ws.Range("J224").Value = "Agree/Accept"
$ws.Range("I274").Value = "sv"
$ws.Range("J274").Value = "Statement-opinion"
$ws.Range("I275").Value = "%"
$ws.Range("J275").Value = "Uninterpretable"
$ws.Range("I276").Value = "%"
$ws.Range("J276").Value = "Uninterpretable"
$ws.Range("I278").Value = "sv"
$ws.Range("J278").Value = "Statement-opinion"
$ws.Range("I283").Value = "sv"
$ws.Range("J283").Value = "Statement-opinion"
$ws.Range("I289").Value = "sv"
$ws.Range("J289").Value = "Statement-opinion"
$ws.Range("I295").Value = "aa"
$ws.Range("J295").Value = "Agree/Accept"
$ws.Range("I306").Value = "sd"
$ws.Range("J306").Value = "Statement-non-opinion"
$ws.Range("I315").Value = "sd"
$ws.Range("J315").Value = "Statement-non-opinion"
$ws.Range("I316").Value = "sd"
$ws.Range("J316").Value = "Statement-non-opinion"
$ws.Range("I321").Value = "sd"
$ws.Range("J321").Value = "Statement-non-opinion"
$ws.Range("I325").Value = "b"
$ws.Range("J325").Value = "Acknowledge (Backchannel)"
$ws.Range("I357").Value = "sd"
$ws.Range("J357").Value = "Statement-non-opinion"
$ws.Range("I362").Value = "aa"
$ws.Range("J362").Value = "Agree/Accept"
$ws.Range("I367").Value = "ba"
$ws.Range("J367").Value = "Appreciation"
$ws.Range("I392").Value = "sd"
$ws.Range("J392").Value = "Statement-non-opinion"
$ws.Range("I408").Value = "sd"
$ws.Range("J408").Value = "Statement-non-opinion"
